# Remove the trailing "Ver no Jupiter ..." / "(c) 2020 ..." boilerplate
# paragraphs (and the blank paragraph that separated them from the last
# bibliography entry), mirroring a Jekyll site rebuild that dropped the
# page-footer scraped from the course listing.

$d = $word.ActiveDocument

# Locate the paragraphs to remove by their text, rather than hard-coded
# indices, so the script is resilient to minor shifts elsewhere in the
# document.
$jupiterIdx = -1
$copyrightIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*Ver no Jupiter*") {
        $jupiterIdx = $i
    }
    if ($t -like "*Creative Commons Attribution*") {
        $copyrightIdx = $i
    }
}

if ($jupiterIdx -gt 0 -and $copyrightIdx -ge $jupiterIdx) {
    # The blank paragraph immediately preceding the "Ver no Jupiter" line
    # is also removed, per the diff.
    $blankIdx = $jupiterIdx - 1

    $rangeStart = $d.Paragraphs.Item($blankIdx).Range.Start
    $rangeEnd = $d.Paragraphs.Item($copyrightIdx).Range.End

    $deleteRange = $d.Range($rangeStart, $rangeEnd)
    $deleteRange.Delete()
}
